$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts boson..diff from F:L to G:M)
$ws.Range("F1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("F1").Value = "pt_max"

# Fill the new column's values (pt_max = 50 for every data row)
$ws.Range("F2:F12").Value = 50

# Update the selection to reflect the newly populated range
$ws.Range("F2:F12").Select()
